# All servers to return the same json payload for valid load tests.
# Update the benchmark numbers on the "Web Servers" sheet with the new
# load-test results (requests / read-MB for 1 Connection and 64 Connections).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# GO / GIN
$ws.Range("D5").Value = 241680
$ws.Range("E5").Value = 127
$ws.Range("F5").Value = 706845
$ws.Range("G5").Value = 372

# CHI
$ws.Range("D6").Value = 270804
$ws.Range("E6").Value = 112
$ws.Range("F6").Value = 939506
$ws.Range("G6").Value = 387

# AIRBORNE
$ws.Range("D7").Value = 262366
$ws.Range("E7").Value = 157
$ws.Range("F7").Value = 743553
$ws.Range("G7").Value = 446

# ECHO
$ws.Range("D8").Value = 249728
$ws.Range("E8").Value = 132
$ws.Range("F8").Value = 849918
$ws.Range("G8").Value = 472

# GORILLA
$ws.Range("D9").Value = 233067
$ws.Range("E9").Value = 96
$ws.Range("F9").Value = 888242
$ws.Range("G9").Value = 366

# FASTHTTP
$ws.Range("D10").Value = 577112
$ws.Range("E10").Value = 248
$ws.Range("F10").Value = 1779269
$ws.Range("G10").Value = 765

# FIBER
$ws.Range("D11").Value = 554514
$ws.Range("E11").Value = 284
$ws.Range("F11").Value = 1571096
$ws.Range("G11").Value = 804

# NODE (JS)
$ws.Range("D13").Value = 33314
$ws.Range("E13").Value = 23
$ws.Range("F13").Value = 52065
$ws.Range("G13").Value = 36

# .NET Core 6.0 webapi (C#)
$ws.Range("D15").Value = 186539
$ws.Range("E15").Value = 116
$ws.Range("F15").Value = 866231
$ws.Range("G15").Value = 541

# Move the active selection (matches the author's saved cursor position).
[void]$ws.Range("J4").Select()
